$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force text-typed assignment (leading apostrophe) without leaving the
# "quote prefix" cell style behind, so the cell style index matches the original.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '66.687.99'
$ws.Cells.Item(2, 5).Value = '  +0.28%  '

Set-TextValue 3 4 '3.256.20'
$ws.Cells.Item(3, 5).Value = '  +2.29%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

Set-TextValue 5 4 '605.57'
$ws.Cells.Item(5, 5).Value = '  +0.23%  '

Set-TextValue 6 4 '159.54'
$ws.Cells.Item(6, 5).Value = '  +2.35%  '

Set-TextValue 7 4 '0.999'
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

Set-TextValue 8 4 '3.256.84'
$ws.Cells.Item(8, 5).Value = '  +2.36%  '

Set-TextValue 9 4 '0.550'
$ws.Cells.Item(9, 5).Value = '  -0.08%  '

$ws.Cells.Item(10, 5).Value = '  +2.11%  '

Set-TextValue 11 4 '5.94'
$ws.Cells.Item(11, 5).Value = '  +4.54%  '

Set-TextValue 12 4 '0.509'
$ws.Cells.Item(12, 5).Value = '  -1.25%  '

Set-TextValue 13 4 '0.0000272'
$ws.Cells.Item(13, 5).Value = '  +2.06%  '

Set-TextValue 14 4 '39.58'
$ws.Cells.Item(14, 5).Value = '  +1.63%  '

Set-TextValue 15 4 '3.786.69'
$ws.Cells.Item(15, 5).Value = '  +2.11%  '

Set-TextValue 16 4 '66.717.57'
$ws.Cells.Item(16, 5).Value = '  +0.26%  '

Set-TextValue 17 4 '7.41'
$ws.Cells.Item(17, 5).Value = '  -0.71%  '

Set-TextValue 18 4 '3.254.49'
$ws.Cells.Item(18, 5).Value = '  +2.22%  '

$ws.Cells.Item(19, 5).Value = '  +1.26%  '

Set-TextValue 20 4 '508.61'
$ws.Cells.Item(20, 5).Value = '  -1.04%  '

Set-TextValue 21 4 '15.41'
$ws.Cells.Item(21, 5).Value = '  -0.95%  '

Set-TextValue 22 4 '0.755'
$ws.Cells.Item(22, 5).Value = '  +2.50%  '

Set-TextValue 23 4 '8.09'
$ws.Cells.Item(23, 5).Value = '  -2.74%  '

Set-TextValue 24 4 '14.89'
$ws.Cells.Item(24, 5).Value = '  -0.70%  '

Set-TextValue 25 4 '86.47'
$ws.Cells.Item(25, 5).Value = '  +2.03%  '

Set-TextValue 26 4 '0.163'
$ws.Cells.Item(26, 5).Value = '  +80.96%  '

$ws.Cells.Item(27, 5).Value = '  +0.14%  '

Set-TextValue 28 4 '3.03'
$ws.Cells.Item(28, 5).Value = '  +0.36%  '

Set-TextValue 29 4 '9.12'
$ws.Cells.Item(29, 5).Value = '  -1.22%  '

Set-TextValue 30 4 '2.40'
$ws.Cells.Item(30, 5).Value = '  -0.75%  '

Set-TextValue 31 4 '6.90'
$ws.Cells.Item(31, 5).Value = '  -2.28%  '

Set-TextValue 32 4 '2.87'
$ws.Cells.Item(32, 5).Value = '  -7.54%  '

Set-TextValue 33 4 '28.21'

$ws.Cells.Item(34, 5).Value = '  -0.03%  '

Set-TextValue 35 4 '1.15'
$ws.Cells.Item(35, 5).Value = '  -4.50%  '

Set-TextValue 36 4 '6.46'
$ws.Cells.Item(36, 5).Value = '  -1.71%  '

Set-TextValue 37 4 '0.0₃0803'
$ws.Cells.Item(37, 5).Value = '  +16.77%  '

# Row 38/39 swap: OKB and dogwifhat swap places (dogwifhat now row 38, OKB now row 39)
$ws.Cells.Item(38, 2).Value = 'dogwifhat'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 38 4 '3.37'
$ws.Cells.Item(38, 5).Value = '  +19.18%  '

$ws.Cells.Item(39, 2).Value = 'OKB'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 39 4 '55.53'
$ws.Cells.Item(39, 5).Value = '  +1.33%  '

Set-TextValue 40 4 '496.80'
$ws.Cells.Item(40, 5).Value = '  -3.43%  '

Set-TextValue 41 4 '0.0429'
$ws.Cells.Item(41, 5).Value = '  +1.22%  '

Set-TextValue 42 4 '0.129'
$ws.Cells.Item(42, 5).Value = '  +1.64%  '

Set-TextValue 43 4 '8.84'
$ws.Cells.Item(43, 5).Value = '  -2.81%  '

Set-TextValue 44 4 '0.297'
$ws.Cells.Item(44, 5).Value = '  -1.48%  '

Set-TextValue 45 4 '2.49'
$ws.Cells.Item(45, 5).Value = '  +1.99%  '

Set-TextValue 46 4 '2.947.92'
$ws.Cells.Item(46, 5).Value = '  +3.27%  '

Set-TextValue 47 4 '28.73'
$ws.Cells.Item(47, 5).Value = '  +0.70%  '

Set-TextValue 48 4 '2.47'
$ws.Cells.Item(48, 5).Value = '  +2.03%  '

Set-TextValue 49 4 '0.120'
$ws.Cells.Item(49, 5).Value = '  +2.00%  '

$ws.Cells.Item(51, 5).Value = '  -2.26%  '
